# Natmi following Dr Hou advice
# Expand the Hp-Itgam LR-pair sheet from 2 data rows (FAPs->M2, FAPs->sCs)
# to 6 data rows covering Sending clusters ECs / FAPs / M2 against
# Target clusters M2 / sCs (Ligand=Hp, Receptor=Itgam throughout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: break every existing reference to the old string pool entries
# (FAPs / Hp / Itgam / M2 / sCs) that live in A2:D3 so they can be
# reintroduced (for the brand new "ECs" cluster) cleanly.
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""

# --- Step 2: introduce the label strings, "ECs" (brand new) first, in the
# same relative order they first appear in the refreshed table.
$ws.Range("A2").Value = "ECs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("D2").Value = "M2"
$ws.Range("B2").Value = "Hp"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D3").Value = "sCs"

# --- Step 3: lay out the full A2:D7 label grid
#   rows 2-3: Sending cluster = ECs   (Target = M2, sCs)
#   rows 4-5: Sending cluster = FAPs  (Target = M2, sCs)
#   rows 6-7: Sending cluster = M2    (Target = M2, sCs)
$labels = @{
    2 = @("ECs",  "Hp", "Itgam", "M2")
    3 = @("ECs",  "Hp", "Itgam", "sCs")
    4 = @("FAPs", "Hp", "Itgam", "M2")
    5 = @("FAPs", "Hp", "Itgam", "sCs")
    6 = @("M2",   "Hp", "Itgam", "M2")
    7 = @("M2",   "Hp", "Itgam", "sCs")
}
foreach ($r in $labels.Keys) {
    $vals = $labels[$r]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
}

# --- Step 4: numeric columns E:T for every data row
$numbers = @{
    2 = @(1, 0.3333333333333333, 0.128414, 0.385242, 0.0205607169378674, 0.0205607169378674, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 5.898224012824, 53.084016115416, 0.02030209810228432, 0.02030209810228432)
    3 = @(1, 0.3333333333333333, 0.128414, 0.385242, 0.0205607169378674, 0.0205607169378674, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 0.07513468896266667, 0.676212200664, 0.0002586188355830755, 0.0002586188355830755)
    4 = @(3, 1, 5.810518333333333, 17.431555, 0.930337990514708, 0.930337990514708, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 266.8847536921267, 2401.96278322914, 0.918635921538578, 0.9186359215385779)
    5 = @(3, 1, 5.810518333333333, 17.431555, 0.930337990514708, 0.930337990514708, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 3.399718782117778, 30.59746903906, 0.01170206897613017, 0.01170206897613017)
    6 = @(2, 0.6666666666666666, 0.306667, 0.9200010000000001, 0.04910129254742459, 0.0491012925474246, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 14.085618883772, 126.770569953948, 0.048483681831679, 0.048483681831679)
    7 = @(2, 0.6666666666666666, 0.306667, 0.9200010000000001, 0.04910129254742459, 0.0491012925474246, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 0.1794300439213334, 1.614870395292, 0.0006176107157455964, 0.0006176107157455965)
}
$numCols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")
foreach ($r in $numbers.Keys) {
    $vals = $numbers[$r]
    for ($i = 0; $i -lt $numCols.Length; $i++) {
        $ws.Range("$($numCols[$i])$r").Value = $vals[$i]
    }
}
